# update lgbm mse space
#
# The diff shows one extra evaluation row added to the "average_mae" sheet: a new
# shared string "ibes_2|fwdepsqcut|tune_mse_ind" is inserted into the workbook's
# string table, and a brand-new row of metrics is inserted at worksheet row 28 -
# every row that used to sit at 28..56 simply moves down one position, to 29..57,
# keeping its own data untouched. Row 28 itself ends up with genuinely new values.
#
# We reproduce that precisely (without leaving behind unused/duplicated cell
# styles, which a plain Rows.Insert() would do) by shifting the existing data
# down manually, from the bottom row upwards, and then writing the brand new
# row's contents into row 28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 28:56 down to 29:57 (bottom-up, so we never overwrite data we still
# need to read).
for ($r = 56; $r -ge 28; $r--) {
    $src = $ws.Range("A" + $r + ":I" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":I" + ($r + 1))
    $dst.Value2 = $src.Value2
}

# Row 57 did not exist before, so its label cell needs the same formatting
# (bold / centered / bordered) as the rest of column A - copy it over from the
# row above, which already has the correct look.
$ws.Cells.Item(56, 1).Copy()
$ws.Cells.Item(57, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Write the brand new row of results into row 28.
$ws.Cells.Item(28, 1).Value = "ibes_2|fwdepsqcut|tune_mse_ind"
$ws.Cells.Item(28, 2).Value = 0.01003309392811347
$ws.Cells.Item(28, 3).Value = 0.009674189396799985
$ws.Cells.Item(28, 4).Value = 0.0003930561638032175
$ws.Cells.Item(28, 5).Value = 0.0003384092866814527
$ws.Cells.Item(28, 6).Value = 0.138487285766222
$ws.Cells.Item(28, 7).Value = 0.2582640092197721
$ws.Cells.Item(28, 8).Value = 0.2582640092197721
$ws.Cells.Item(28, 9).Value = 14156
